# Actualizado 12-mar con nuevas fases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Alevín A): jornada 12 -> jornada 13
$ws.Range("A6").Value = "https://www.rfebm.com/competiciones/competicion.php?seleccion=0&id=1026965&jornada=13&id_ambito=0"

# Row 8 (Infantil F): new results link (id 1026259 -> id 1028401)
$ws.Range("A8").Value = "https://www.rfebm.com/competiciones/resultados_completos.php?seleccion=0&id=1028401"

# New empty cell H13, formatted like the hyperlink-style cells (underlined text, same as A13)
$ws.Range("H13").Font.Underline = 2

# Leave the selection on the newly touched cell, as in the saved workbook
$ws.Range("H13").Select()
